$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wiki_Cats_to_Topics")

# Sort the D2:E25 range by column D ascending (blanks sort last)
$rngDE = $ws.Range("D2:E25")
$keyD = $ws.Range("D2:D25")
$rngDE.Sort($keyD, 1, $null, $null, 1, $null, 1, 0, $false, $null, $null, 1)

# Sort the G2:H18 range by column G ascending (blanks sort last)
$rngGH = $ws.Range("G2:H18")
$keyG = $ws.Range("G2:G18")
$rngGH.Sort($keyG, 1, $null, $null, 1, $null, 1, 0, $false, $null, $null, 1)

# Set column D width
$ws.Columns.Item(4).ColumnWidth = 19

# Update selection / view
$ws.Range("H18").Select()
$ws.Application.ActiveWindow.ScrollColumn = 5

$wb.Windows.Item(1).WindowState = -4143
